$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New consolidated values for rows 2-4
$ws.Range("A2").Value = "('Argothian Enchantress', ['{1}{G}', 'Creature — Human Druid', 'Shroud (This creature can’t be the target of spells or abilities.)', 'Whenever you cast an enchantment spell, draw a card.', '0/1'])"
$ws.Range("A3").Value = "('Intuition', ['{2}{U}', 'Instant', 'Search your library for three cards and reveal them. Target opponent chooses one. Put that card into your hand and the rest into your graveyard. Then shuffle your library.'])"
$ws.Range("A4").Value = "('Living Death', ['{3}{B}{B}', 'Sorcery', 'Each player exiles all creature cards from their graveyard, then sacrifices all creatures they control, then puts all cards they exiled this way onto the battlefield.'])"

# Remove old rows 5-15 which are no longer needed
$ws.Range("A5:A15").EntireRow.Delete()
